$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (SignUp_Positive_01): Condition flips from "No" to "Yes",
# and the auto-generated email/firstname/lastname values are refreshed.
$ws.Range("C2").Value = "Yes"
$ws.Range("I2").Value = "prabhaAutohzLW5080@mailinator.com"
$ws.Range("K2").Value = "PrabhaAutoZjVh"
$ws.Range("L2").Value = "automationlebs"

# Row 5 (SignUp_Positive_04): Condition flips from "Yes" to "No".
$ws.Range("C5").Value = "No"

# Update the saved selection to C2 (matches the authored workbook state).
$ws.Range("C2").Select()
